$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (exhibitions) - update a handful of "want to go" counters
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 13134
$wsExpo.Range("F8").Value = 33
$wsExpo.Range("F10").Value = 13091
$wsExpo.Range("F13").Value = 8787
$wsExpo.Range("F14").Value = 7843
$wsExpo.Range("F19").Value = 997

# ---------------------------------------------------------------------------
# Sheet "演出" (performances) - update one counter
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 21

# ---------------------------------------------------------------------------
# Sheet "全部类型" (all types) - same counters plus a brand-new event row
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 13134
$wsAll.Range("F9").Value = 33
$wsAll.Range("F11").Value = 13091
$wsAll.Range("F14").Value = 8787
$wsAll.Range("F15").Value = 7843

# Insert a brand new row 20 (shifting the existing rows 20-29 down to 21-30)
# for a newly scraped event: 昆山·心动次元动漫游戏嘉年华
$wsAll.Rows.Item(20).Insert(-4121)

# Index column - match the look/feel (bold, centered, bordered) of the other
# index cells in column A.
$wsAll.Range("A20").Value = 19
$wsAll.Range("A20").Font.Bold = $true
$wsAll.Range("A20").HorizontalAlignment = -4108
$wsAll.Range("A20").VerticalAlignment = -4160
$wsAll.Range("A20").Borders.LineStyle = 1

# Date column must stay plain text (like every other date cell in the sheet)
# instead of being auto-converted to a date serial number.
$wsAll.Range("B20").NumberFormat = "@"
$wsAll.Range("B20").Value = "2024-10-03"
$wsAll.Range("B20").Style = $wsAll.Range("C19").Style

$wsAll.Range("C20").Value = "昆山·心动次元动漫游戏嘉年华"
$wsAll.Range("D20").Value = "娄苑路413号 叁加叁球馆"
$wsAll.Range("E20").Value = "2024.10.03 10:00-10.03 17:00"
$wsAll.Range("F20").Value = 0
$wsAll.Range("G20").Value = 45
$wsAll.Range("H20").Value = "https://show.bilibili.com/platform/detail.html?id=91643"
$wsAll.Range("I20").Value = "//i0.hdslb.com/bfs/openplatform/202408/HQiYu4i81724642509360.jpeg"

# The index numbers in column A are plain literal values (not a formula), so
# every row that got pushed down by the insert needs to be renumbered by +1.
for ($r = 21; $r -le 30; $r++) {
    $cell = $wsAll.Range("A$r")
    $cell.Value = $cell.Value2 + 1
}
